$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-ParaIndexByContains($doc, $substring) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.Contains($substring)) { return $i }
    }
    return -1
}

# Inserts a brand-new (empty) ListParagraph-style paragraph immediately after
# the paragraph at $idx, copying the list level ($ilvl, 0-based as in OOXML)
# onto it, and returns the 1-based paragraph index of the new paragraph.
function Insert-ListParagraphAfter($doc, $idx, $ilvl) {
    $anchor = $doc.Paragraphs.Item($idx)
    $rng = $anchor.Range
    $rng.Collapse(0) | Out-Null
    $rng.InsertParagraphAfter() | Out-Null
    $newIdx = $idx + 1
    $newP = $doc.Paragraphs.Item($newIdx)
    $newP.Range.ListFormat.ListLevelNumber = $ilvl + 1
    return $newIdx
}

# Sets the (non-hyperlink) text of paragraph $idx to $text, keeping the
# trailing paragraph mark intact.
function Set-ParagraphText($doc, $idx, $text) {
    $r = $doc.Paragraphs.Item($idx).Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

# Turns the (already-set) text currently occupying paragraph $idx into a
# hyperlink pointing at $address.
function Add-ParagraphHyperlink($doc, $idx, $address) {
    $r = $doc.Paragraphs.Item($idx).Range
    $r.MoveEnd(1, -1) | Out-Null
    $doc.Hyperlinks.Add($r, $address, "", "", $address) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Insert a new sub-bullet right after "Disease by country CSV from ESPEN"
#    containing a hyperlink to the ESPEN country page, plus a note, and a
#    further sub-sub-bullet with navigation instructions.
#    (Doing this before the GBD insertion -- which sits earlier in the
#    document -- keeps paragraph indices below it stable.)
# ---------------------------------------------------------------------------

$espenAnchorIdx = Get-ParaIndexByContains $d "Disease by country CSV from ESPEN"

$espenLinkIdx = Insert-ListParagraphAfter $d $espenAnchorIdx 1
Set-ParagraphText $d $espenLinkIdx "http://espen.afro.who.int/countries/nigeria"
Add-ParagraphHyperlink $d $espenLinkIdx "http://espen.afro.who.int/countries/nigeria"
$tail = $d.Paragraphs.Item($espenLinkIdx).Range
$tail.MoveEnd(1, -1) | Out-Null
$tail.Collapse(0) | Out-Null
$tail.InsertAfter(" (replace Nigeria with country name)") | Out-Null

$navigateIdx = Insert-ListParagraphAfter $d $espenLinkIdx 2
Set-ParagraphText $d $navigateIdx "Navigate to the disease you are looking at, and then the data option."

# ---------------------------------------------------------------------------
# 2) Insert a new sub-bullet right after the "YLD Rate pulled..." paragraph
#    containing a hyperlink to the GBD results tool.
# ---------------------------------------------------------------------------

$gbdAnchorIdx = Get-ParaIndexByContains $d "YLD Rate pulled from the GBD tool"

$gbdLinkIdx = Insert-ListParagraphAfter $d $gbdAnchorIdx 1
Set-ParagraphText $d $gbdLinkIdx "http://ghdx.healthdata.org/gbd-results-tool"
Add-ParagraphHyperlink $d $gbdLinkIdx "http://ghdx.healthdata.org/gbd-results-tool"

# ---------------------------------------------------------------------------
# 3) Re-word the two "future ideas" bullets at the end of the document.
# ---------------------------------------------------------------------------

$futureMetricIdx = Get-ParaIndexByContains $d "Future idea to create a new metric"
Set-ParagraphText $d $futureMetricIdx "Drugs Used vs Procured"

$mapIdx = Get-ParaIndexByContains $d "Makes the map even easier to read"
Set-ParagraphText $d $mapIdx "DrugPopTreat / Drug Total Procured"

# The last bullet used to carry a (now pointless) "_GoBack" bookmark; drop it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

Write-Host "Edit complete."
